# Redo report on laptop
# Updates the n_ratings (column G) counts for a set of rows in the
# pairs report sheet, reflecting re-collected ratings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11  = 2
    14  = 1
    19  = 1
    27  = 2
    28  = 1
    29  = 1
    30  = 1
    32  = 2
    41  = 1
    42  = 1
    46  = 2
    47  = 2
    48  = 2
    51  = 1
    54  = 1
    56  = 2
    63  = 2
    69  = 1
    71  = 1
    75  = 1
    76  = 2
    88  = 2
    94  = 2
    99  = 2
    103 = 3
    104 = 2
    109 = 2
    115 = 3
    127 = 2
    128 = 2
    129 = 2
    132 = 1
    148 = 1
    155 = 3
    157 = 2
    165 = 2
    173 = 3
    180 = 3
    183 = 1
    185 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
